# Auto-generated Excel COM-interop edit script
# Applies the numeric cell updates for the scheduled market-data refresh
# across the ALC/ARM/CRP/CUL/GSM/LTW/WVR sheets (see commit message).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 33357
$ws.Range("J3").Value = 33357
$ws.Range("L3").Value = 33357
$ws.Range("N3").Value = -33585

$ws.Range("H32").Value = 1950
$ws.Range("I32").Value = 900
$ws.Range("J32").Value = 3000
$ws.Range("K32").Value = 900
$ws.Range("L32").Value = 3000
$ws.Range("M32").Value = -574
$ws.Range("N32").Value = -3652

$ws.Range("H74").Value = 5257.5713
$ws.Range("I74").Value = 5634.3335
$ws.Range("J74").Value = 4975
$ws.Range("K74").Value = 5634.3335
$ws.Range("L74").Value = 4975
$ws.Range("M74").Value = -4698.3335
$ws.Range("N74").Value = -6847

$ws.Range("H76").Value = 3875.6667
$ws.Range("I76").Value = 3500.5
$ws.Range("J76").Value = 4250.8335
$ws.Range("K76").Value = 3500.5
$ws.Range("L76").Value = 4250.8335
$ws.Range("M76").Value = -3185.5
$ws.Range("N76").Value = -4880.8335

$ws.Range("H77").Value = 5257.5713
$ws.Range("I77").Value = 5634.3335
$ws.Range("J77").Value = 4975
$ws.Range("K77").Value = 28171.6675
$ws.Range("L77").Value = 24875
$ws.Range("M77").Value = -23491.6675
$ws.Range("N77").Value = -34235

$ws.Range("H79").Value = 3875.6667
$ws.Range("I79").Value = 3500.5
$ws.Range("J79").Value = 4250.8335
$ws.Range("K79").Value = 3500.5
$ws.Range("L79").Value = 4250.8335
$ws.Range("M79").Value = -2408.5
$ws.Range("N79").Value = -6434.8335

$ws.Range("H98").Value = 2204.862
$ws.Range("I98").Value = 1277.44
$ws.Range("J98").Value = 8001.25
$ws.Range("K98").Value = 1277.44
$ws.Range("L98").Value = 8001.25
$ws.Range("M98").Value = 220.5599999999999
$ws.Range("N98").Value = -10997.25

$ws.Range("H99").Value = 1307.8334
$ws.Range("I99").Value = 370.5
$ws.Range("J99").Value = 2245.1667
$ws.Range("K99").Value = 1111.5
$ws.Range("L99").Value = 6735.500100000001
$ws.Range("M99").Value = 386.5
$ws.Range("N99").Value = -9731.500100000001

$ws.Range("H102").Value = 33357
$ws.Range("J102").Value = 33357
$ws.Range("L102").Value = 33357
$ws.Range("N102").Value = -39847

$ws.Range("H122").Value = 2204.862
$ws.Range("I122").Value = 1277.44
$ws.Range("J122").Value = 8001.25
$ws.Range("K122").Value = 3832.32
$ws.Range("L122").Value = 24003.75
$ws.Range("M122").Value = -1382.32
$ws.Range("N122").Value = -28903.75

$ws.Range("H129").Value = 3379794.5
$ws.Range("I129").Value = 27778750
$ws.Range("J129").Value = 1477.5231
$ws.Range("K129").Value = 83336250
$ws.Range("L129").Value = 4432.5693
$ws.Range("M129").Value = -83331250
$ws.Range("N129").Value = -14432.5693

$ws.Range("H137").Value = 3230
$ws.Range("I137").Value = 3200
$ws.Range("J137").Value = 3260
$ws.Range("K137").Value = 9600
$ws.Range("L137").Value = 9780
$ws.Range("M137").Value = -7050
$ws.Range("N137").Value = -14880

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4481.918
$ws.Range("I32").Value = 3217.238
$ws.Range("K32").Value = 3217.238
$ws.Range("M32").Value = -2930.238

$ws.Range("H98").Value = 24999.75
$ws.Range("J98").Value = 24999.75
$ws.Range("L98").Value = 24999.75
$ws.Range("N98").Value = -30989.75

$ws.Range("H101").Value = 29344.727
$ws.Range("J101").Value = 29344.727
$ws.Range("L101").Value = 29344.727
$ws.Range("N101").Value = -35834.727

$ws.Range("H104").Value = 141750
$ws.Range("J104").Value = 141750
$ws.Range("L104").Value = 141750
$ws.Range("N104").Value = -148738

$ws.Range("H122").Value = 4271.2
$ws.Range("I122").Value = 2542.4
$ws.Range("J122").Value = 6000
$ws.Range("K122").Value = 7627.200000000001
$ws.Range("L122").Value = 18000
$ws.Range("M122").Value = -5177.200000000001
$ws.Range("N122").Value = -22900

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 3796
$ws.Range("J9").Value = 3796
$ws.Range("L9").Value = 3796
$ws.Range("N9").Value = -4132

$ws.Range("H58").Value = 15627713
$ws.Range("I58").Value = 1855.2916
$ws.Range("J58").Value = 62505290
$ws.Range("K58").Value = 1855.2916
$ws.Range("L58").Value = 62505290
$ws.Range("M58").Value = -1652.2916
$ws.Range("N58").Value = -62505696

$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").Value = $null

$ws.Range("H132").Value = 3249.1086
$ws.Range("I132").Value = 1963.2759
$ws.Range("J132").Value = 5442.5884
$ws.Range("K132").Value = 5889.8277
$ws.Range("L132").Value = 16327.7652
$ws.Range("M132").Value = -3359.8277
$ws.Range("N132").Value = -21387.7652

$ws.Range("H136").Value = 15627713
$ws.Range("I136").Value = 1855.2916
$ws.Range("J136").Value = 62505290
$ws.Range("K136").Value = 5565.8748
$ws.Range("L136").Value = 187515870
$ws.Range("M136").Value = -3015.8748
$ws.Range("N136").Value = -187520970

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 1807
$ws.Range("J92").Value = 1933.3125
$ws.Range("L92").Value = 5799.9375
$ws.Range("N92").Value = -8295.9375

$ws.Range("H107").Value = 1677.25
$ws.Range("J107").Value = 2999.5
$ws.Range("L107").Value = 8998.5
$ws.Range("N107").Value = -12838.5

$ws.Range("H119").Value = 3265.2144
$ws.Range("I119").Value = 1682
$ws.Range("J119").Value = 3898.5
$ws.Range("K119").Value = 5046
$ws.Range("L119").Value = 11695.5
$ws.Range("M119").Value = -208
$ws.Range("N119").Value = -21371.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3819
$ws.Range("I70").Value = 3779.8572
$ws.Range("J70").Value = 3873.8
$ws.Range("K70").Value = 3779.8572
$ws.Range("L70").Value = 3873.8
$ws.Range("M70").Value = -3509.8572
$ws.Range("N70").Value = -4413.8

$ws.Range("H73").Value = 3819
$ws.Range("I73").Value = 3779.8572
$ws.Range("J73").Value = 3873.8
$ws.Range("K73").Value = 3779.8572
$ws.Range("L73").Value = 3873.8
$ws.Range("M73").Value = -2843.8572
$ws.Range("N73").Value = -5745.8

$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").Value = $null

$ws.Range("H126").Value = 2429.7026
$ws.Range("I126").Value = 1395.6086
$ws.Range("J126").Value = 4128.5713
$ws.Range("K126").Value = 4186.825800000001
$ws.Range("L126").Value = 12385.7139
$ws.Range("M126").Value = -1716.825800000001
$ws.Range("N126").Value = -17325.7139

$ws.Range("H132").Value = 2945.138
$ws.Range("I132").Value = 2506.4119
$ws.Range("J132").Value = 3566.6667
$ws.Range("K132").Value = 7519.2357
$ws.Range("L132").Value = 10700.0001
$ws.Range("M132").Value = -4989.2357
$ws.Range("N132").Value = -15760.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H131").Value = 20783.334
$ws.Range("J131").Value = 20783.334
$ws.Range("L131").Value = 20783.334
$ws.Range("N131").Value = -30863.334

$ws.Range("H133").Value = 29490
$ws.Range("J133").Value = 29490
$ws.Range("L133").Value = 29490
$ws.Range("N133").Value = -34550

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 20000
$ws.Range("J27").Value = 20000
$ws.Range("L27").Value = 20000
$ws.Range("N27").Value = -20138

$ws.Range("H115").Value = 38248.75
$ws.Range("J115").Value = 38248.75
$ws.Range("L115").Value = 38248.75
$ws.Range("N115").Value = -41382.75

$ws.Range("H126").Value = 3227480.2
$ws.Range("I126").Value = 1369.3478
$ws.Range("J126").Value = 12502549
$ws.Range("K126").Value = 4108.0434
$ws.Range("L126").Value = 37507647
$ws.Range("M126").Value = -3185.4546
$ws.Range("N126").Value = -37512587
